$wb = $excel.ActiveWorkbook

# --- Sheet: no_of_dependents ---
$lastIdx = $wb.Worksheets.Count
$ws1 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIdx))
$ws1.Name = "no_of_dependents"

$ws1.Range("A1").Value = "no_of_dependents"
$ws1.Range("B1").Value = "count"

$data1 = @(
  @(4, 752),
  @(3, 727),
  @(0, 712),
  @(2, 708),
  @(1, 697),
  @(5, 673)
)
for ($i = 0; $i -lt $data1.Count; $i++) {
  $row = $i + 2
  $ws1.Cells.Item($row, 1).Value = $data1[$i][0]
  $ws1.Cells.Item($row, 2).Value = $data1[$i][1]
}

$headerRange1 = $ws1.Range("A1:B1")
$headerRange1.Font.Bold = $true
$headerRange1.HorizontalAlignment = -4108
$headerRange1.VerticalAlignment = -4160
$headerRange1.Borders.LineStyle = 1
$headerRange1.Borders.Weight = 2

$colARange1 = $ws1.Range("A2:A7")
$colARange1.Font.Bold = $true
$colARange1.HorizontalAlignment = -4108
$colARange1.VerticalAlignment = -4160
$colARange1.Borders.LineStyle = 1
$colARange1.Borders.Weight = 2

# --- Sheet: loan_status ---
$lastIdx2 = $wb.Worksheets.Count
$ws2 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIdx2))
$ws2.Name = "loan_status"

$ws2.Range("A1").Value = "loan_status"
$ws2.Range("B1").Value = "count"

$ws2.Range("A2").Value = " Approved"
$ws2.Range("B2").Value = 2656
$ws2.Range("A3").Value = " Rejected"
$ws2.Range("B3").Value = 1613

$headerRange2 = $ws2.Range("A1:B1")
$headerRange2.Font.Bold = $true
$headerRange2.HorizontalAlignment = -4108
$headerRange2.VerticalAlignment = -4160
$headerRange2.Borders.LineStyle = 1
$headerRange2.Borders.Weight = 2

$colARange2 = $ws2.Range("A2:A3")
$colARange2.Font.Bold = $true
$colARange2.HorizontalAlignment = -4108
$colARange2.VerticalAlignment = -4160
$colARange2.Borders.LineStyle = 1
$colARange2.Borders.Weight = 2
